$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.886.57"
$ws.Range("E2").Value = "  +0.45%  "

$ws.Range("D3").Value = "3.073.14"
$ws.Range("E3").Value = "  +2.03%  "

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "0.996"
$c.Style = "Normal"
$ws.Range("E4").Value = "  -0.47%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "522.98"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +2.48%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "141.13"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +1.23%  "

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.997"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -0.28%  "

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.436"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -0.26%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "7.37"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -2.03%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.111"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +0.94%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.374"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +2.52%  "

$ws.Range("D12").Value = "3.559.26"
$ws.Range("E12").Value = "  +0.94%  "

$ws.Range("E13").Value = "  -3.12%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "27.13"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +2.58%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.0000170"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +3.34%  "

$ws.Range("D16").Value = "57.773.46"
$ws.Range("E16").Value = "  +0.30%  "

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "6.28"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +0.67%  "

$ws.Range("D18").Value = "3.054.60"
$ws.Range("E18").Value = "  +1.38%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "13.47"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +5.08%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "8.27"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +3.79%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "333.71"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +0.76%  "

$ws.Range("E22").Value = "  +0.21%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "0.513"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +2.86%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "65.10"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +0.87%  "

$ws.Range("D25").Value = "3.161.52"
$ws.Range("E25").Value = "  +0.86%  "

$ws.Range("B26").Value = "Kaspa"
$ws.Range("C26").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "0.166"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -2.11%  "

$ws.Range("B27").Value = "Binance-PegBSC-USD"
$ws.Range("C27").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -0.02%  "

$ws.Range("D28").Value = "0.0₃0913"
$ws.Range("E28").Value = "  -1.04%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "6.77"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -0.17%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "7.35"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -0.45%  "

$ws.Range("B31").Value = "Fetch.AI"
$ws.Range("C31").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "1.24"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +4.67%  "

$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "1.83"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +0.73%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "21.06"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +2.16%  "

$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "4.69"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -0.54%  "

$ws.Range("B35").Value = "Monero"
$ws.Range("C35").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "154.36"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +0.25%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "5.94"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +1.38%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "25.87"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +5.87%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "1.29"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +1.04%  "

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.0681"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +0.82%  "

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "37.36"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +0.02%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "3.93"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +2.11%  "

$ws.Range("B42").Value = "Mantle"
$ws.Range("C42").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.667"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +2.51%  "

$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.998"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -0.34%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "1.41"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +0.23%  "

$ws.Range("D45").Value = "2.208.73"
$ws.Range("E45").Value = "  -1.33%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "6.21"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +3.20%  "

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.970"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -1.50%  "

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "0.0246"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +2.54%  "

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "20.18"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +4.38%  "

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "1.80"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -2.94%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.186"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +1.85%  "
